# Add 2022-Q3 data
# 1) Insert a new "2022-Q3" sheet before the existing "2022-Q2" sheet and
#    populate it with the fund-holdings detail for that quarter.
# 2) Insert a new summary row on the "总计" sheet for 2022-Q3, pushing the
#    existing quarters down and renumbering the index column.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Part 1: "总计" summary sheet - insert a 2022-Q3 row at the top of the
# data (right under the header row).
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()

# Reuse the existing header/index cell formatting (bold, centered, boxed)
# instead of re-deriving it, so no new style entries are created.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 9
$total.Range("D2").Value = 2.53

# Renumber the index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# ---------------------------------------------------------------------
# Part 2: new "2022-Q3" sheet with the fund holdings detail.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($q2)
$q3.Name = "2022-Q3"

# Borrow the header/index-column style (s="2") from the "总计" sheet so it
# matches the workbook's existing look without minting new style entries.
$total.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$q3.Range("A2:A10").PasteSpecial(-4122)

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Columns B (fund code, has leading zeros) and D:G (numeric-looking
# percentages/sizes) are text in the source data - force text formatting
# first so Excel doesn't silently coerce them to numbers.
$q3.Range("B2:B10").NumberFormat = "@"
$q3.Range("D2:G10").NumberFormat = "@"

$data = @(
  @(0, "014254", "信澳智远三年持有期混合A", "36.34", "86.96", "1.45", "0.5269", 10),
  @(1, "012608", "信澳领先智选混合",         "33.22", "92.57", "1.43", "0.4750", 10),
  @(2, "007484", "信澳核心科技混合",         "23.33", "94.48", "1.80", "0.4199", 7),
  @(3, "006257", "信澳先进智造股票",         "25.85", "92.65", "1.61", "0.4162", 9),
  @(4, "011188", "信澳星奕混合A",            "23.94", "93.29", "1.45", "0.3471", 10),
  @(5, "009511", "信澳研究优选混合A",         "10.05", "89.66", "1.43", "0.1437", 10),
  @(6, "011223", "信澳星奕混合C",            "7.06",  "93.29", "1.45", "0.1024", 10),
  @(7, "014255", "信澳智远三年持有期混合C",   "6.26",  "86.96", "1.45", "0.0908", 10),
  @(8, "014954", "信澳研究优选混合C",         "0.22",  "89.66", "1.43", "0.0031", 10)
)

$r = 2
foreach ($row in $data) {
  $q3.Range("A$r").Value = $row[0]
  $q3.Range("B$r").Value = $row[1]
  $q3.Range("C$r").Value = $row[2]
  $q3.Range("D$r").Value = $row[3]
  $q3.Range("E$r").Value = $row[4]
  $q3.Range("F$r").Value = $row[5]
  $q3.Range("G$r").Value = $row[6]
  $q3.Range("H$r").Value = $row[7]
  $r = $r + 1
}
